# Financials update: refresh yearly figures on the ASWRF income/balance/cash-flow sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ASWRF")

$ws.Range("G17").Value = 1500
$ws.Range("J17").Value = 2300
$ws.Range("G18").Value = -1500
$ws.Range("J18").Value = -2300
$ws.Range("G20").Value = -5500
$ws.Range("G21").Value = -7000
$ws.Range("F23").Value = -1800
$ws.Range("G23").Value = -7000
$ws.Range("J23").Value = -2300
$ws.Range("F26").Value = -1800
$ws.Range("G26").Value = -7000
$ws.Range("I26").Value = -800
$ws.Range("F27").Value = -1800
$ws.Range("G27").Value = -7000
$ws.Range("I27").Value = -800
$ws.Range("G32").Value = 5500
$ws.Range("F33").Value = -1800
$ws.Range("G33").Value = -7000
$ws.Range("I33").Value = -800
$ws.Range("F35").Value = -1800
$ws.Range("G35").Value = -7000
$ws.Range("I35").Value = -800
$ws.Range("H48").Value = 8800
$ws.Range("I48").Value = 8900
$ws.Range("J48").Value = 9400
$ws.Range("G54").Value = 2400
$ws.Range("H54").Value = 9000
$ws.Range("I54").Value = 9500
$ws.Range("J54").Value = 10300
$ws.Range("I57").Value = 100
$ws.Range("I60").Value = 100
$ws.Range("I66").Value = 100
$ws.Range("D72").Value = -26400
$ws.Range("E72").Value = -26700
$ws.Range("F72").Value = -26400
$ws.Range("G72").Value = -24600
$ws.Range("H72").Value = -16800
$ws.Range("I72").Value = -15600
$ws.Range("J72").Value = -14500
$ws.Range("H76").Value = 9000
$ws.Range("I76").Value = 9400
$ws.Range("J76").Value = 10000
$ws.Range("F81").Value = -1800
$ws.Range("G81").Value = -7000
$ws.Range("I81").Value = -800
$ws.Range("D94").Value = -700
$ws.Range("I94").Value = -200
$ws.Range("J94").Value = -600
$ws.Range("I102").Value = -500
